$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row written as new content (previously the sheet was empty)
$ws.Range("A1").Value = "USN"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Age"
$ws.Range("D1").Value = "Phone"
$ws.Range("E1").Value = "Vaccine_Dose"

# Column E widened to fit the longer "Vaccine_Dose" header
$ws.Range("E1").EntireColumn.AutoFit() | Out-Null

# Cursor left on the next empty cell after the header row
$ws.Range("F1").Select() | Out-Null
